$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as plain text, preserving the cell's original
# (unstyled / General) appearance. A bare `.Value = "62.00"` assignment
# would get silently re-interpreted as the number 62, dropping the
# trailing zero / thousands-style formatting used throughout this sheet
# (prices like "59.119.29", "62.00", "0.998" are text, not numbers).
# Flipping to a Text number format for the write and then back to the
# Normal style keeps the stored value a string without leaving the cell
# tagged with a custom style.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "59.119.29"
Set-TextValue $ws.Range("E2") "  -2.97%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "2.544.84"
Set-TextValue $ws.Range("E3") "  +0.40%  "

# Row 4 - TetherUSD
Set-TextValue $ws.Range("D4") "0.998"
Set-TextValue $ws.Range("E4") "  -0.60%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "538.96"
Set-TextValue $ws.Range("E5") "  -0.10%  "

# Row 6 - Solana
Set-TextValue $ws.Range("D6") "144.13"
Set-TextValue $ws.Range("E6") "  -3.23%  "

# Row 7 - USDC
Set-TextValue $ws.Range("E7") "  -0.58%  "

# Row 8 - XRP
Set-TextValue $ws.Range("D8") "0.572"
Set-TextValue $ws.Range("E8") "  -0.41%  "

# Row 9 - LidoStakedEther
Set-TextValue $ws.Range("D9") "2.543.93"
Set-TextValue $ws.Range("E9") "  +0.11%  "

# Row 10 - Dogecoin
Set-TextValue $ws.Range("D10") "0.0997"
Set-TextValue $ws.Range("E10") "  -0.98%  "

# Row 11 - TRON
Set-TextValue $ws.Range("E11") "  -1.45%  "

# Row 12 - Toncoin
Set-TextValue $ws.Range("D12") "5.53"
Set-TextValue $ws.Range("E12") "  +2.20%  "

# Row 13 - Cardano
Set-TextValue $ws.Range("D13") "0.352"
Set-TextValue $ws.Range("E13") "  -1.01%  "

# Row 14 - WrappedliquidstakedEther2.0
Set-TextValue $ws.Range("D14") "2.951.23"
Set-TextValue $ws.Range("E14") "  -1.37%  "

# Row 15 - Avalanche
Set-TextValue $ws.Range("E15") "  -4.26%  "

# Row 16 - WrappedBTC
Set-TextValue $ws.Range("D16") "58.951.90"
Set-TextValue $ws.Range("E16") "  -2.96%  "

# Row 17 - ShibaInu
Set-TextValue $ws.Range("D17") "0.0000139"
Set-TextValue $ws.Range("E17") "  -0.25%  "

# Row 18 - WrappedEther
Set-TextValue $ws.Range("D18") "2.535.25"
Set-TextValue $ws.Range("E18") "  -0.52%  "

# Row 19 - Chainlink
Set-TextValue $ws.Range("D19") "11.31"
Set-TextValue $ws.Range("E19") "  -0.37%  "

# Row 20 - Polkadot
Set-TextValue $ws.Range("E20") "  -2.93%  "

# Row 21 - BitcoinCash
Set-TextValue $ws.Range("D21") "324.13"
Set-TextValue $ws.Range("E21") "  -1.88%  "

# Row 22 - Dai
Set-TextValue $ws.Range("E22") "  -0.12%  "

# Row 23 - Uniswap
Set-TextValue $ws.Range("D23") "5.79"
Set-TextValue $ws.Range("E23") "  -1.25%  "

# Row 24 - Litecoin
Set-TextValue $ws.Range("D24") "62.00"
Set-TextValue $ws.Range("E24") "  -0.29%  "

# Row 25 - Polygon
Set-TextValue $ws.Range("E25") "  -7.71%  "

# Row 26 - Kaspa
Set-TextValue $ws.Range("E26") "  -0.24%  "

# Row 27 - WrappedeETH
Set-TextValue $ws.Range("D27") "2.626.19"
Set-TextValue $ws.Range("E27") "  -1.94%  "

# Row 28 - Binance-PegBSC-USD
Set-TextValue $ws.Range("D28") "0.992"
Set-TextValue $ws.Range("E28") "  +0.05%  "

# Row 29 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D29") "7.78"
Set-TextValue $ws.Range("E29") "  -1.61%  "

# Row 30 - Aptos
Set-TextValue $ws.Range("E30") "  -3.50%  "

# Row 31 - PEPE
Set-TextValue $ws.Range("E31") "  -1.84%  "

# Row 32 - PancakeSwap
Set-TextValue $ws.Range("E32") "  -3.76%  "

# Row 33 - Fetch.AI
Set-TextValue $ws.Range("D33") "1.21"
Set-TextValue $ws.Range("E33") "  -6.00%  "

# Row 34 - USDe
Set-TextValue $ws.Range("E34") "  -0.42%  "

# Row 35 - Monero
Set-TextValue $ws.Range("D35") "158.02"
Set-TextValue $ws.Range("E35") "  -1.92%  "

# Row 36 - ImmutableX
Set-TextValue $ws.Range("D36") "1.44"
Set-TextValue $ws.Range("E36") "  +5.44%  "

# Row 37 - EthereumClassic
Set-TextValue $ws.Range("E37") "  -0.45%  "

# Row 38 - NEARProtocol
Set-TextValue $ws.Range("D38") "4.39"
Set-TextValue $ws.Range("E38") "  -5.50%  "

# Row 39 - Stacks
Set-TextValue $ws.Range("E39") "  -6.14%  "

# Row 40 - RenderToken
Set-TextValue $ws.Range("D40") "5.70"
Set-TextValue $ws.Range("E40") "  -2.43%  "

# Row 41 - Bittensor
Set-TextValue $ws.Range("D41") "303.98"
Set-TextValue $ws.Range("E41") "  -3.99%  "

# Row 42 - OKB
Set-TextValue $ws.Range("D42") "36.87"
Set-TextValue $ws.Range("E42") "  -0.33%  "

# Row 43 - SuiNetwork
Set-TextValue $ws.Range("E43") "  -5.23%  "

# Row 44 - Filecoin
Set-TextValue $ws.Range("E44") "  -4.01%  "

# Row 45 - FirstDigitalUSD
Set-TextValue $ws.Range("D45") "0.998"
Set-TextValue $ws.Range("E45") "  -0.97%  "

# Row 46 - Mantle
Set-TextValue $ws.Range("E46") "  +1.38%  "

# Row 47 - WhiteBITCoin
Set-TextValue $ws.Range("E47") "  -1.03%  "

# Row 48 - Aave
Set-TextValue $ws.Range("D48") "126.39"
Set-TextValue $ws.Range("E48") "  +5.05%  "

# Row 49 - Stellar
Set-TextValue $ws.Range("D49") "0.0932"
Set-TextValue $ws.Range("E49") "  -1.78%  "

# Row 50 - EnergySwap
Set-TextValue $ws.Range("D50") "18.74"
Set-TextValue $ws.Range("E50") "  -1.25%  "

# Row 51 - Hedera
Set-TextValue $ws.Range("D51") "0.0517"
Set-TextValue $ws.Range("E51") "  -1.79%  "
